# Actualización automática 2026-01-07 16:30:09
#
# Inserts a new salesperson row ("VERA CABRERA JORGE ENRIQUE") above the
# existing "VIEJO RIVAS MAYRA ANABELLE" row on both worksheets, shifting
# the trailing summary row down, and updates the summary row's totals /
# "0 de N" counters to reflect the extra row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Push the existing row 22 (and the totals row below it) down by one,
# duplicating row 22's formatting into the newly opened row.
$ws1.Rows.Item(22).Insert()

$ws1.Range("A22").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Range("B22").Value = "VERA CABRERA JORGE ENRIQUE"

$cols1 = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols1) {
    $ws1.Range($col + "22").Value = 0
}

# Update the "0 de 21" -> "0 de 22" counters on what is now row 24.
foreach ($col in $cols1) {
    $ws1.Range($col + "24").Value = "0 de 22"
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(22).Insert()

$ws2.Range("A22").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Range("B22").Value = "VERA CABRERA JORGE ENRIQUE"
$ws2.Range("C22").Value = 1115.07
$ws2.Range("D22").Value = 0
$ws2.Range("E22").Value = 388.8
$ws2.Range("F22").Value = 0
$ws2.Range("G22").Value = 1000

# Update the totals row, now at row 24.
$ws2.Range("C24").Value = 18337.91
$ws2.Range("D24").Value = 3750.24
$ws2.Range("E24").Value = 16490.86
$ws2.Range("F24").Value = 0
$ws2.Range("G24").Value = 1000

$wb.Save()
